$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Propagate the existing header format (bold white-on-green, centered,
#        wrapped - currently on B2) onto the new header cells A1 and B1 ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Write the new header text ---
$ws.Range("A1").Value = "DOCUMENTO "
$ws.Range("B1").Value = "FECHA DE CESE (OPCIONAL): 2020-12-10"

# --- 3. Clear the old "DNI " text/format out of B2 - it becomes a plain,
#        unstyled input cell below the new headers ---
$ws.Range("B2").ClearContents() | Out-Null
$ws.Range("B2").ClearFormats() | Out-Null

# --- 4. Column widths to fit the new, longer headers ---
$ws.Columns.Item(1).ColumnWidth = 20.75
$ws.Columns.Item(2).ColumnWidth = 40.75

# --- 5. Selection matches the authored file (active cell B2) ---
$ws.Range("B2").Select() | Out-Null

# --- 6. Page setup (portrait / letter, matching the authored file) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
